$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 719, pushing the existing rows 719-765 down to 720-766.
$ws.Rows.Item(719).Insert()

# Populate the newly inserted row 719 with the new data record.
$ws.Cells.Item(719, 1).Value = 4
$ws.Cells.Item(719, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(719, 3).Value = "Los Lagos"
$ws.Cells.Item(719, 4).Value = 44931
$ws.Cells.Item(719, 5).Value = 10
$ws.Cells.Item(719, 6).Value = 100112004
$ws.Cells.Item(719, 7).Value = "Cebolla"
$ws.Cells.Item(719, 8).Value = "Sin especificar"
$ws.Cells.Item(719, 9).Value = "1a nueva(o)"
$ws.Cells.Item(719, 10).Value = 400
$ws.Cells.Item(719, 11).Value = 17000
$ws.Cells.Item(719, 12).Value = 17000
$ws.Cells.Item(719, 13).Value = 17000
$ws.Cells.Item(719, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(719, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(719, 16).Value = 944
$ws.Cells.Item(719, 17).Value = 18
$ws.Cells.Item(719, 18).Value = "Hortaliza"
